$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112 (ALC)
$ws.Range("H112").Value = 1265.3448
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 1285.5358
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 3856.6074
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -6072.607400000001

# Row 135 (ALC)
$ws.Range("H135").Value = 563.9535
$ws.Range("I135").Value = 541.90247
$ws.Range("J135").Value = 1016
$ws.Range("K135").Value = 4877.12223
$ws.Range("L135").Value = 9144
$ws.Range("M135").Value = -2342.12223
$ws.Range("N135").Value = -14214

# Row 137 (ALC)
$ws.Range("H137").Value = 1631.4783
$ws.Range("I137").Value = 1343.6428
$ws.Range("J137").Value = 2079.2222
$ws.Range("K137").Value = 4030.9284
$ws.Range("L137").Value = 6237.6666
$ws.Range("M137").Value = -1480.9284
$ws.Range("N137").Value = -11337.6666

# Row 138 (ALC)
$ws.Range("H138").Value = 2195.1843
$ws.Range("I138").Value = 1513.1936
$ws.Range("J138").Value = 2665
$ws.Range("K138").Value = 4539.5808
$ws.Range("L138").Value = 7995
$ws.Range("M138").Value = 600.4192000000003
$ws.Range("N138").Value = -18275

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 4092.953
$ws.Range("I32").Value = 3401.4324
$ws.Range("K32").Value = 3401.4324
$ws.Range("M32").Value = -3114.4324

# Row 61 (ARM)
$ws.Range("H61").Value = 1821.0857
$ws.Range("I61").Value = 1447.909
$ws.Range("J61").Value = 2452.6155
$ws.Range("K61").Value = 1447.909
$ws.Range("L61").Value = 2452.6155
$ws.Range("M61").Value = -1235.909
$ws.Range("N61").Value = -2876.6155

# Row 74 (ARM)
$ws.Range("H74").Value = 596.2973
$ws.Range("I74").Value = 607.6857
$ws.Range("K74").Value = 607.6857
$ws.Range("M74").Value = 266.3143

# Row 77 (ARM)
$ws.Range("H77").Value = 596.2973
$ws.Range("I77").Value = 607.6857
$ws.Range("K77").Value = 3038.4285
$ws.Range("M77").Value = 1329.5715

# Row 95 (ARM)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 132 (ARM)
$ws.Range("H132").Value = 4501.625
$ws.Range("I132").Value = 5227.926
$ws.Range("J132").Value = 2993.1538
$ws.Range("K132").Value = 15683.778
$ws.Range("L132").Value = 8979.4614
$ws.Range("M132").Value = -13153.778
$ws.Range("N132").Value = -14039.4614

# Row 136 (ARM)
$ws.Range("H136").Value = 1821.0857
$ws.Range("I136").Value = 1447.909
$ws.Range("J136").Value = 2452.6155
$ws.Range("K136").Value = 4343.727000000001
$ws.Range("L136").Value = 7357.8465
$ws.Range("M136").Value = -1793.727000000001
$ws.Range("N136").Value = -12457.8465

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (BSM)
$ws.Range("H134").Value = 39522.555
$ws.Range("I134").Value = 61037.94
$ws.Range("J134").Value = 2946.4
$ws.Range("K134").Value = 183113.82
$ws.Range("L134").Value = 8839.200000000001
$ws.Range("M134").Value = -180578.82
$ws.Range("N134").Value = -13909.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 3970649
$ws.Range("I31").Value = 1446.7878
$ws.Range("J31").Value = 18524390
$ws.Range("K31").Value = 1446.7878
$ws.Range("L31").Value = 18524390
$ws.Range("M31").Value = -1151.7878
$ws.Range("N31").Value = -18524980

# Row 34 (CRP)
$ws.Range("H34").Value = 3970649
$ws.Range("I34").Value = 1446.7878
$ws.Range("J34").Value = 18524390
$ws.Range("K34").Value = 1446.7878
$ws.Range("L34").Value = 18524390
$ws.Range("M34").Value = -1244.7878
$ws.Range("N34").Value = -18524794

# Row 58 (CRP)
$ws.Range("H58").Value = 855.8611
$ws.Range("I58").Value = 676.2759
$ws.Range("J58").Value = 1599.8572
$ws.Range("K58").Value = 676.2759
$ws.Range("L58").Value = 1599.8572
$ws.Range("M58").Value = -473.2759
$ws.Range("N58").Value = -2005.8572

# Row 132 (CRP)
$ws.Range("H132").Value = 1712.5294
$ws.Range("I132").Value = 1374.6154
$ws.Range("J132").Value = 2810.75
$ws.Range("K132").Value = 4123.8462
$ws.Range("L132").Value = 8432.25
$ws.Range("M132").Value = -1593.8462
$ws.Range("N132").Value = -13492.25

# Row 134 (CRP)
$ws.Range("H134").Value = 1532.5333
$ws.Range("I134").Value = 1548.8
$ws.Range("K134").Value = 4646.4
$ws.Range("M134").Value = -2111.4

# Row 136 (CRP)
$ws.Range("H136").Value = 855.8611
$ws.Range("I136").Value = 676.2759
$ws.Range("J136").Value = 1599.8572
$ws.Range("K136").Value = 2028.8277
$ws.Range("L136").Value = 4799.571599999999
$ws.Range("M136").Value = 521.1723000000002
$ws.Range("N136").Value = -9899.571599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 32 (CUL)
$ws.Range("H32").Value = 200
$ws.Range("I32").Value = 200
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -317
$ws.Range("N32").ClearContents()

# Row 121 (CUL)
$ws.Range("H121").Value = 1038.5
$ws.Range("I121").Value = 182.11111
$ws.Range("J121").Value = 1202.4894
$ws.Range("K121").Value = 546.3333299999999
$ws.Range("L121").Value = 3607.4682
$ws.Range("M121").Value = 763.6666700000001
$ws.Range("N121").Value = -6227.468199999999

# Row 131 (CUL)
$ws.Range("H131").Value = 1545475.2
$ws.Range("I131").Value = 4812.8
$ws.Range("J131").Value = 2364976.5
$ws.Range("K131").Value = 14438.4
$ws.Range("L131").Value = 7094929.5
$ws.Range("M131").Value = -9398.400000000001
$ws.Range("N131").Value = -7105009.5

# Row 137 (CUL)
$ws.Range("H137").Value = 55608384
$ws.Range("I137").Value = 30304054
$ws.Range("J137").Value = 75490360
$ws.Range("K137").Value = 90912162
$ws.Range("L137").Value = 226471080
$ws.Range("M137").Value = -90907062
$ws.Range("N137").Value = -226481280

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (LTW)
$ws.Range("H132").Value = 5501.744
$ws.Range("I132").Value = 7008.0347
$ws.Range("J132").Value = 2381.5715
$ws.Range("K132").Value = 21024.1041
$ws.Range("L132").Value = 7144.7145
$ws.Range("M132").Value = -18494.1041
$ws.Range("N132").Value = -12204.7145

# Row 136 (LTW)
$ws.Range("H136").Value = 3414.8635
$ws.Range("I136").Value = 3866.1292
$ws.Range("J136").Value = 2338.7693
$ws.Range("K136").Value = 11598.3876
$ws.Range("L136").Value = 7016.3079
$ws.Range("M136").Value = -12116.3079

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws.Range("H96").Value = 25000710
$ws.Range("I96").Value = 83333700
$ws.Range("J96").Value = 857.1429000000001
$ws.Range("K96").Value = 83333700
$ws.Range("L96").Value = 857.1429000000001
$ws.Range("M96").Value = -83332327
$ws.Range("N96").Value = -3603.1429

# Row 132 (WVR)
$ws.Range("H132").Value = 1183.02
$ws.Range("I132").Value = 959.8333
$ws.Range("J132").Value = 2354.75
$ws.Range("K132").Value = 2879.4999
$ws.Range("L132").Value = 7064.25
$ws.Range("M132").Value = -349.4998999999998
$ws.Range("N132").Value = -12124.25

# Row 136 (WVR)
$ws.Range("H136").Value = 3637.0435
$ws.Range("I136").Value = 3948.5264
$ws.Range("K136").Value = 11845.5792
$ws.Range("M136").Value = -9295.5792
